$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.175.36'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '1.908.46'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''325.88'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').Value = '''0.4626'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''0.3895'
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('D9').Value = '''0.07882'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('D10').Value = '''0.9925'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').Value = '''21.99'
$ws.Range('E11').Value = '  -1.96%  '
$ws.Range('D12').Value = '1.889.44'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').Value = '''5.769'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').Value = '''7.054'
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('D15').Value = '''0.07042'
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('D16').Value = '''88.20'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '''1.003'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '''0.000009950'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').Value = '''17.12'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '29.183.88'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').Value = '''5.328'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').Value = '''11.14'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.107.49'
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '''2.095'
$ws.Range('E25').Value = '  +1.44%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''156.56'
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''19.48'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''5.915'
$ws.Range('E28').Value = '  -3.43%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '''119.04'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = '''1.884'
$ws.Range('E30').Value = '  -5.92%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '''0.09358'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '''0.8977'
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''5.236'
$ws.Range('E33').Value = '  -2.34%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '''1.325'
$ws.Range('E34').Value = '  -2.74%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '''3.147'
$ws.Range('E35').Value = '  -3.97%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.05797'
$ws.Range('E36').Value = '  -0.80%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '''1.175'
$ws.Range('E37').Value = '  -3.11%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.02094'
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('B39').Value = 'Frax'
$ws.Range('C39').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D39').Value = '''1.001'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.5712'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''7.671'
$ws.Range('E41').Value = '  -4.10%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = '''0.1816'
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '''9.745'
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''11.87'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '''0.5363'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '''2.174'
$ws.Range('E46').Value = '  -5.64%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '''0.07014'
$ws.Range('E47').Value = '  -0.94%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''1.845'
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = '''2.549'
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''113.33'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '''0.2977'
$ws.Range('E51').Value = '  +0.80%  '
